# Stock summary update (raj) - apply quantity/rate/value corrections
# to the PATRIKA / CARD rows called out in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B = Quantity, C = Rate, D = factor, E = Value

# Row 9  - 1670 PATRIKA (06)
$ws.Cells.Item(9, 2).Value = 62
$ws.Cells.Item(9, 3).Value = 277
$ws.Cells.Item(9, 5).Value = 277

# Row 51 - 1889 PATRIKA (Eco 55)
$ws.Cells.Item(51, 2).Value = 55
$ws.Cells.Item(51, 3).Value = 7
$ws.Cells.Item(51, 5).Value = 7.72

# Row 52 - 1890 PATRIKA (Eco 56)
$ws.Cells.Item(52, 2).Value = 364
$ws.Cells.Item(52, 3).Value = 13
$ws.Cells.Item(52, 5).Value = 14.63

# Row 81 - 1914 PATRIKA (YELLOW)
$ws.Cells.Item(81, 2).Value = 256
$ws.Cells.Item(81, 3).Value = 20.5

# Row 82 - 1915 PATRIKA (RED)
$ws.Cells.Item(82, 2).Value = 307
$ws.Cells.Item(82, 3).Value = 618.5
$ws.Cells.Item(82, 5).Value = 711.28

# Row 84 - 1917 PATRIKA (RED)
$ws.Cells.Item(84, 2).Value = 457
$ws.Cells.Item(84, 3).Value = 1988.8
$ws.Cells.Item(84, 5).Value = 2784.32

# Row 86 - 1919 PATRIKA (BLUE)
$ws.Cells.Item(86, 2).Value = 174
$ws.Cells.Item(86, 3).Value = 600
$ws.Cells.Item(86, 5).Value = 840

# Row 90 - 1922 PATRIKA
$ws.Cells.Item(90, 2).Value = 250
$ws.Cells.Item(90, 3).Value = 23
$ws.Cells.Item(90, 5).Value = 28.86

# Row 98 - 1930 PATRIKA (S.S.63)
$ws.Cells.Item(98, 2).Value = 35
$ws.Cells.Item(98, 3).Value = 55.5
$ws.Cells.Item(98, 5).Value = 86.03

# Row 163 - 4214 PATRIKA
$ws.Cells.Item(163, 2).Value = 28
$ws.Cells.Item(163, 3).Value = 15.3
$ws.Cells.Item(163, 5).Value = 42.84

# Row 250 - 5046 PATRIKA
$ws.Cells.Item(250, 3).Value = 2.5
$ws.Cells.Item(250, 5).Value = 9.5

# Row 258 - 5856 PATRIKA
$ws.Cells.Item(258, 2).Value = 245
$ws.Cells.Item(258, 3).Value = 23.75
$ws.Cells.Item(258, 5).Value = 101.65

# Row 272 - 5869 PATRIKA
$ws.Cells.Item(272, 2).Value = 76
$ws.Cells.Item(272, 3).Value = 8
$ws.Cells.Item(272, 5).Value = 44.4

# Row 383 - 6531 PATRIKA {M}
$ws.Cells.Item(383, 2).Value = 53
$ws.Cells.Item(383, 3).Value = 0.49
$ws.Cells.Item(383, 5).Value = 12.25

# Row 413 - 7258 PATRIKA
$ws.Cells.Item(413, 2).Value = 22
$ws.Cells.Item(413, 3).Value = 29
$ws.Cells.Item(413, 5).Value = 65.25

# Row 414 - 7259 PATRIKA
$ws.Cells.Item(414, 2).Value = 31
$ws.Cells.Item(414, 3).Value = 26
$ws.Cells.Item(414, 5).Value = 58.5

# Row 418 - 7263 PATRIKA : quantity bumped, rate/factor/value now blank (no sale)
$ws.Cells.Item(418, 2).Value = 21
$ws.Cells.Item(418, 3).ClearContents()
$ws.Cells.Item(418, 4).ClearContents()
$ws.Cells.Item(418, 5).ClearContents()
$ws.Cells.Item(418, 3).NumberFormat = $ws.Cells.Item(12, 3).NumberFormat
$ws.Cells.Item(418, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(418, 5).NumberFormat = $ws.Cells.Item(12, 5).NumberFormat

# Row 433 - 7278 PATRIKA
$ws.Cells.Item(433, 2).Value = 66
$ws.Cells.Item(433, 3).Value = 21
$ws.Cells.Item(433, 5).Value = 66.150000000000006

# Row 441 - 7286 PATRIKA (HALKA D/F)
$ws.Cells.Item(441, 2).Value = 73
$ws.Cells.Item(441, 3).Value = 19
$ws.Cells.Item(441, 5).Value = 66.5

# Row 602 - 9097 CARD
$ws.Cells.Item(602, 2).Value = 18
$ws.Cells.Item(602, 3).Value = 11.5
$ws.Cells.Item(602, 5).Value = 49.22

# Row 613 - 9108 CARDS
$ws.Cells.Item(613, 2).Value = 24
$ws.Cells.Item(613, 3).Value = 589

# Row 617 - 9201 CARD (GOLDEN)
$ws.Cells.Item(617, 2).Value = 155
$ws.Cells.Item(617, 3).Value = 101.5
$ws.Cells.Item(617, 5).Value = 81.2

# Row 618 - 9202 CARD (YELLOW)
$ws.Cells.Item(618, 2).Value = 105
$ws.Cells.Item(618, 3).Value = 240.5
$ws.Cells.Item(618, 5).Value = 195.94

# Row 639 - Grand Total
$ws.Cells.Item(639, 3).Value = 22697.63
$ws.Cells.Item(639, 5).Value = 38546.82
